$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label tweak ---
$ws.Range("C2").Value = "Short Notes"

# --- Wording/title fixes on existing rows ---
$ws.Range("B5").Value = "Remove Duplicates from Sorted Array"
$ws.Range("B7").Value = "Remove Element"

# --- Row 11 used to hold the "9. Arrays / Second Maximum Element" entry;
#     that entry is dropped and replaced with a new LeetCode problem. ---
$ws.Range("A11").ClearContents()
$ws.Range("B11").Value = "Palindrome Number"
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = "https://leetcode.com/problems/palindrome-number/"

# --- New rows for the 100-days-of-coding additions ---
$ws.Range("B12").Value = "Roman to Integer"
$ws.Range("D12").Value = "https://leetcode.com/problems/roman-to-integer/"

$ws.Range("B13").Value = "Merge Two Sorted Lists"
$ws.Range("D13").Value = "https://leetcode.com/problems/merge-two-sorted-lists/"

$ws.Range("B14").Value = "Maximum Subarray"
$ws.Range("D14").Value = "https://leetcode.com/problems/maximum-subarray/"

# Match the "Good" highlight styling used by the other problem-name cells
# (B6:B10 use this same look) for the new / repurposed rows.
$ws.Range("B11").Style = "Good"
$ws.Range("B12").Style = "Good"
$ws.Range("B13").Style = "Good"
$ws.Range("B14").Style = "Good"

# --- Selection moves to C27 ---
$ws.Range("C27").Select()

Write-Output "done"
